$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2..19): System name, BO, BI, UI, AO, AI, Pressure
$data = @(
    @("UC600_test",  4, 0,  8, 6, 0, 1),
    @("UC600_test2", 4, 0,  9, 5, 0, 0),
    @("UC600_test",  4, 0,  7, 7, 0, 1),
    @("UC600_test3", 5, 0,  6, 6, 0, 0),
    @("UC600_test",  0, 14, 0, 0, 0, 1),
    @("UC600_test4", 4, 14, 0, 0, 0, 2),
    @("UC600_test",  4, 0, 14, 0, 0, 1),
    @("UC600_test5", 4, 0, 14, 0, 0, 0),
    @("s500test1",   9, 0,  0, 0, 0, 2),
    @("s500test2",  10, 0,  0, 0, 0, 0),
    @("s500test3",   9, 3,  0, 0, 0, 0),
    @("s500test4",   9, 5,  0, 0, 0, 0),
    @("s500test5",   9, 7,  0, 0, 0, 2),
    @("s500test6",   9, 5,  0, 2, 0, 0),
    @("s500test7",   9, 3,  2, 3, 3, 0),
    @("s500test8",   9, 3,  0, 2, 7, 0),
    @("s500test9",   9, 3,  3, 3, 3, 2),
    @("s500test10",  9, 4,  2, 3, 5, 1)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}

$ws.Columns.Item(1).ColumnWidth = 12

$ws.Range("H19").Select()
